$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting existing rows 107..238 down to 108..239.
$ws.Rows.Item(107).EntireRow.Insert()

# Populate the newly inserted row 107 with its data (a new Perejil price
# observation at Vega Modelo de Temuco, dated 2021-12-16).
$ws.Cells.Item(107, 1).Value = 10
$ws.Cells.Item(107, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(107, 3).Value = "La Araucanía"
$ws.Cells.Item(107, 4).Value = 44546
$ws.Cells.Item(107, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(107, 5).Value = 9
$ws.Cells.Item(107, 6).Value = 100112044
$ws.Cells.Item(107, 7).Value = "Perejil"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 65
$ws.Cells.Item(107, 11).Value = 4500
$ws.Cells.Item(107, 12).Value = 4500
$ws.Cells.Item(107, 13).Value = 4500
$ws.Cells.Item(107, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(107, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(107, 16).Value = 1500
$ws.Cells.Item(107, 17).Value = 3
$ws.Cells.Item(107, 18).Value = "Hortaliza"
